# "basic setup done" - turn the blank Sheet1 into a watchlist of NSE
# scripts with RSI/EMA columns, highlight the header row, and rename
# the tab.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Watchlist"

# Header row
$ws.Range("A1").Value = "Scripts"
$ws.Range("B1").Value = "Close"
$ws.Range("C1").Value = "RSI(14)"
$ws.Range("D1").Value = "EMA(50)"
$ws.Range("E1").Value = "EMA(200)"

# SBIN
$ws.Range("A2").Value = "SBIN"
$ws.Range("C2").Value = 47.78
$ws.Range("D2").Value = 466.85
$ws.Range("E2").Value = 471.57

# HDFCBANK
$ws.Range("A3").Value = "HDFCBANK"
$ws.Range("C3").Value = 39.65
$ws.Range("D3").Value = 1382.92
$ws.Range("E3").Value = 1396.07

# ICICIBANK
$ws.Range("A4").Value = "ICICIBANK"
$ws.Range("C4").Value = 33.11
$ws.Range("D4").Value = 748.14
$ws.Range("E4").Value = 752.96

# ZOMATO
$ws.Range("A5").Value = "ZOMATO"
$ws.Range("C5").Value = 62.39
$ws.Range("D5").Value = 76.64
$ws.Range("E5").Value = 73.97

# Highlight the header row yellow
$ws.Range("A1:E1").Interior.Color = 65535

# Widen the Scripts column so the longer names aren't clipped
$ws.Columns("A").ColumnWidth = 11.14

# Leave the selection on EMA(50) like the saved workbook
$ws.Range("D1").Select() | Out-Null
